$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45183 -> 45184) for every data row from row 2 through row 20.
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}
